$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = "z495256"

$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "z495256"

$ws.Range("B5").Value = "z495257"

[void]$ws.Range("C13").Select()
